$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 262
$ws.Range("I28").Value = 212.8
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 212.8
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 272.2
$ws.Range("N28").Value = -1970

$ws.Range("H107").Value = 455.93332
$ws.Range("I107").Value = 468.2
$ws.Range("J107").Value = 449.8
$ws.Range("K107").Value = 468.2
$ws.Range("L107").Value = 449.8
$ws.Range("M107").Value = 1451.8
$ws.Range("N107").Value = -4289.8

$ws.Range("H129").Value = 182667.73
$ws.Range("J129").Value = 182667.73
$ws.Range("L129").Value = 548003.1900000001
$ws.Range("N129").Value = -558003.1900000001

$ws.Range("H132").Value = 4313.3
$ws.Range("I132").Value = 4641.25
$ws.Range("K132").Value = 13923.75
$ws.Range("M132").Value = -11393.75

$ws.Range("H137").Value = 94463.09
$ws.Range("I137").Value = 6250
$ws.Range("J137").Value = 144870.58
$ws.Range("K137").Value = 18750
$ws.Range("L137").Value = 434611.74
$ws.Range("M137").Value = -16200
$ws.Range("N137").Value = -439711.74

$ws.Range("H138").Value = 1560.324
$ws.Range("J138").Value = 1964.4722
$ws.Range("L138").Value = 5893.4166
$ws.Range("N138").Value = -16173.4166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18890.441
$ws.Range("I32").Value = 19537.475
$ws.Range("J32").Value = 450
$ws.Range("K32").Value = 19537.475
$ws.Range("L32").Value = 450
$ws.Range("M32").Value = -19250.475
$ws.Range("N32").Value = -1024

$ws.Range("H45").Value = 3737.7273
$ws.Range("I45").Value = 4880.2
$ws.Range("J45").Value = 2785.6667
$ws.Range("K45").Value = 4880.2
$ws.Range("L45").Value = 2785.6667
$ws.Range("M45").Value = -4503.2
$ws.Range("N45").Value = -3539.6667

$ws.Range("H61").Value = 2632.8572
$ws.Range("I61").Value = 2281.2632
$ws.Range("K61").Value = 2281.2632
$ws.Range("M61").Value = -2069.2632

$ws.Range("H74").Value = 55558424
$ws.Range("I74").Value = 55558424
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 55558424
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -55557550
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 55558424
$ws.Range("I77").Value = 55558424
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 277792120
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -277787752
$ws.Range("N77").ClearContents()

$ws.Range("H110").Value = 637.05884
$ws.Range("I110").Value = 547.53845
$ws.Range("J110").Value = 928
$ws.Range("K110").Value = 547.53845
$ws.Range("L110").Value = 928
$ws.Range("M110").Value = 1497.46155
$ws.Range("N110").Value = -5018

$ws.Range("H122").Value = 2112.5789
$ws.Range("I122").Value = 2189.0625
$ws.Range("K122").Value = 6567.1875
$ws.Range("M122").Value = -4117.1875

$ws.Range("H132").Value = 15550
$ws.Range("I132").Value = 2031.4783
$ws.Range("K132").Value = 6094.4349
$ws.Range("M132").Value = -3564.4349

$ws.Range("H136").Value = 2632.8572
$ws.Range("I136").Value = 2281.2632
$ws.Range("K136").Value = 6843.7896
$ws.Range("M136").Value = -4293.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2833.3333
$ws.Range("J99").Value = 2750
$ws.Range("L99").Value = 2750
$ws.Range("N99").Value = -5746

$ws.Range("H107").Value = 1592.4
$ws.Range("I107").Value = 801.375
$ws.Range("K107").Value = 801.375
$ws.Range("M107").Value = 1118.625

$ws.Range("H134").Value = 20002.666
$ws.Range("I134").Value = 26694.096
$ws.Range("J134").Value = 1266.6666
$ws.Range("K134").Value = 80082.288
$ws.Range("L134").Value = 3799.9998
$ws.Range("M134").Value = -77547.288
$ws.Range("N134").Value = -8869.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13509.16
$ws.Range("I31").Value = 15425.19
$ws.Range("K31").Value = 15425.19
$ws.Range("M31").Value = -15130.19

$ws.Range("H34").Value = 13509.16
$ws.Range("I34").Value = 15425.19
$ws.Range("K34").Value = 15425.19
$ws.Range("M34").Value = -15223.19

$ws.Range("H59").Value = 20163.316
$ws.Range("I59").Value = 7551.5
$ws.Range("J59").Value = 21647.059
$ws.Range("K59").Value = 7551.5
$ws.Range("L59").Value = 21647.059
$ws.Range("M59").Value = -6406.5
$ws.Range("N59").Value = -23937.059

$ws.Range("H132").Value = 19726.793
$ws.Range("I132").Value = 29219
$ws.Range("K132").Value = 87657
$ws.Range("M132").Value = -85127

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 33666.668
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340

$ws.Range("H129").Value = 358631.44
$ws.Range("I129").Value = 672.75
$ws.Range("J129").Value = 501814.9
$ws.Range("K129").Value = 2018.25
$ws.Range("L129").Value = 1505444.7
$ws.Range("M129").Value = 2981.75
$ws.Range("N129").Value = -1515444.7

$ws.Range("H131").Value = 736.05
$ws.Range("J131").Value = 736.05
$ws.Range("L131").Value = 2208.15
$ws.Range("N131").Value = -12288.15

$ws.Range("H137").Value = 2695
$ws.Range("I137").Value = 1196.5
$ws.Range("K137").Value = 3589.5
$ws.Range("M137").Value = 1510.5

$ws.Range("H140").Value = 1723.6111
$ws.Range("I140").Value = 1530.625
$ws.Range("J140").Value = 3267.5
$ws.Range("K140").Value = 4591.875
$ws.Range("L140").Value = 9802.5
$ws.Range("M140").Value = 588.125
$ws.Range("N140").Value = -20162.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31252708
$ws.Range("I102").Value = 38464650
$ws.Range("J102").Value = 978
$ws.Range("K102").Value = 38464650
$ws.Range("L102").Value = 978
$ws.Range("M102").Value = -38463028
$ws.Range("N102").Value = -4222

$ws.Range("H122").Value = 148149500
$ws.Range("I122").Value = 83334920
$ws.Range("K122").Value = 250004760
$ws.Range("M122").Value = -250002310

$ws.Range("H132").Value = 106530.734
$ws.Range("I132").Value = 107671.3
$ws.Range("K132").Value = 323013.9
$ws.Range("M132").Value = -320483.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5972.222
$ws.Range("I7").Value = 6342.5
$ws.Range("K7").Value = 6342.5
$ws.Range("M7").Value = -6230.5

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H93").Value = 2266.7144
$ws.Range("I93").Value = 2043.8334
$ws.Range("J93").Value = 3604
$ws.Range("K93").Value = 2043.8334
$ws.Range("L93").Value = 3604
$ws.Range("M93").Value = -795.8334
$ws.Range("N93").Value = -6100

$ws.Range("H122").Value = 1228765.5
$ws.Range("I122").Value = 1636461.9
$ws.Range("K122").Value = 4909385.699999999
$ws.Range("M122").Value = -4906935.699999999

$ws.Range("H126").Value = 5972.222
$ws.Range("I126").Value = 6342.5
$ws.Range("K126").Value = 19027.5
$ws.Range("M126").Value = -16557.5

$ws.Range("H132").Value = 1989.7778
$ws.Range("I132").Value = 1272.4286
$ws.Range("K132").Value = 3817.2858
$ws.Range("M132").Value = -1287.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4049.7
$ws.Range("I62").Value = 2500.3333
$ws.Range("K62").Value = 2500.3333
$ws.Range("M62").Value = -1876.3333

$ws.Range("H63").Value = 39998
$ws.Range("J63").Value = 39998
$ws.Range("L63").Value = 39998
$ws.Range("N63").Value = -41246

$ws.Range("H65").Value = 4049.7
$ws.Range("I65").Value = 2500.3333
$ws.Range("K65").Value = 12501.6665
$ws.Range("M65").Value = -9381.666499999999

$ws.Range("H66").Value = 39998
$ws.Range("J66").Value = 39998
$ws.Range("L66").Value = 119994
$ws.Range("N66").Value = -126234

$ws.Range("H107").Value = 3247702.2
$ws.Range("I107").Value = 1185.25
$ws.Range("K107").Value = 3555.75
$ws.Range("M107").Value = -1635.75

$ws.Range("H122").Value = 1534.4445
$ws.Range("I122").Value = 1626.2916
$ws.Range("J122").Value = 1350.75
$ws.Range("K122").Value = 4878.8748
$ws.Range("L122").Value = 4052.25
$ws.Range("M122").Value = -2428.8748
